$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (row 5)
$ws.Range("A5").Value = 64
$ws.Range("B5").Value = "Hot Wheels Program"
$ws.Range("C5").Value = "disability centers"
$ws.Range("D5").Value = "310 Thrift Rd"
$ws.Range("E5").Value = "Madison"
$ws.Range("G5").Value = -78.262205399999999
$ws.Range("H5").Value = 38.3771542

# Select the newly added row, matching the saved selection state
$ws.Rows.Item(5).Select()
